$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "DAILY MEAL DATA" to "MONTHLY MEAL DATA"
$ws.Name = "MONTHLY MEAL DATA"

# Row 2
$ws.Range("A2").Value = 28615
$ws.Range("B2").Value = "John Smith"
$ws.Range("C2").Value = "Diet"
$ws.Range("D2").Value = "1:20 PM"
$ws.Range("E2").Value = "1 July 2024"
$ws.Range("F2").Value = "Karachi"

# Row 3
$ws.Range("A3").Value = 28615
$ws.Range("B3").Value = "John Smith"
$ws.Range("C3").Value = "Diet"
$ws.Range("D3").Value = "1:20 PM"
$ws.Range("E3").Value = "18 July 2024"
$ws.Range("F3").Value = "Karachi"

# Row 4
$ws.Range("A4").Value = 31489
$ws.Range("B4").Value = "Sara Malik"
$ws.Range("C4").Value = "Normal"
$ws.Range("D4").Value = "1:50 PM"
$ws.Range("E4").Value = "1 July 2024"
$ws.Range("F4").Value = "Karachi"

# Row 5
$ws.Range("A5").Value = 28615
$ws.Range("B5").Value = "John Smith"
$ws.Range("C5").Value = "Diet"
$ws.Range("D5").Value = "2:20 PM"
$ws.Range("E5").Value = "10 July 2024"
$ws.Range("F5").Value = "Karachi"

# Row 6
$ws.Range("A6").Value = 28615
$ws.Range("B6").Value = "John Smith"
$ws.Range("C6").Value = "Diet"
$ws.Range("D6").Value = "2:20 PM"
$ws.Range("E6").Value = "18 July 2024"
$ws.Range("F6").Value = "Karachi"
